$p = $ppt.ActivePresentation
Write-Host "HasTitleMaster=$($p.HasTitleMaster)"
$tm = $p.TitleMaster
Write-Host "tm=$tm"
$tcs = $tm.Theme.ThemeColorScheme
Write-Host "Count=$($tcs.Count)"
Write-Host "dk2=$($tcs.Item(3).RGB)"
